$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new "Paid" column data (adds Paid/Y/N to sharedStrings, extends the table to col E) ----
$ws.Range("E1").Value = "Paid"
$ws.Range("E2").Value = "Y"
$ws.Range("E3").Value = "N"
$ws.Range("E4").Value = "N"
$ws.Range("E5").Value = "Y"
$ws.Range("E6").Value = "Y"
$ws.Range("E7").Value = "Y"

# ---- column widths ----
$ws.Columns(1).ColumnWidth = 9.83
$ws.Columns(2).ColumnWidth = 19
$ws.Columns(3).ColumnWidth = 21.83
$ws.Columns(4).ColumnWidth = 17.67

# ---- row heights ----
$ws.Rows(1).RowHeight = 19
$ws.Rows(2).RowHeight = 19
$ws.Rows(3).RowHeight = 17
$ws.Rows(4).RowHeight = 22
$ws.Rows(5).RowHeight = 17
$ws.Rows(6).RowHeight = 17
$ws.Rows(7).RowHeight = 17

# ---- per-cell formatting (font, fill, border, alignment) ----
$c = $ws.Range("A1")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 14
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 4
$c.Interior.Color = 1137349
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("B1")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 14
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 4
$c.Interior.Color = 1137349
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("C1")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 14
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 4
$c.Interior.Color = 1137349
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("D1")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 14
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 4
$c.Interior.Color = 1137349
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("E1")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 14
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 4
$c.Interior.Color = 1137349
$c.HorizontalAlignment = -4108

$c = $ws.Range("A2")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("B2")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 14
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("C2")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("D2")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4152

$c = $ws.Range("E2")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("A3")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("B3")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("C3")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("D3")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4152

$c = $ws.Range("E3")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("A4")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("B4")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 16
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("C4")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Interior.Pattern = -4142
$c.Borders.LineStyle = 1

$c = $ws.Range("D4")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $true
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.Color = 65535
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4152

$c = $ws.Range("E4")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("A5")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("B5")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("C5")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("D5")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $true
$c.Font.Italic = $true
$c.Font.Underline = $false
$c.Font.Color = 11534591
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4152

$c = $ws.Range("E5")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("A6")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("B6")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $true
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("C6")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("D6")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4152

$c = $ws.Range("E6")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

$c = $ws.Range("A7")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("B7")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $true
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("C7")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1

$c = $ws.Range("D7")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4152

$c = $ws.Range("E7")
$c.Font.Name = "Courier New"
$c.Font.Family = 1
$c.Font.Size = 12
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Font.Underline = $false
$c.Font.ThemeColor = 1
$c.Borders.LineStyle = 1
$c.HorizontalAlignment = -4108

# ---- selection / window ----
[void]$ws.Range("A1:E7").Select()

